# Fix in loar area coordinators
# Updates the two sample user rows (row 2: "Pedro Paredes", row 3: "Mary James")
# with corrected data, fixes the hyperlink display/target for the email
# column, adds a missing hyperlink for the new row, widens column D, and
# moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: Pedro Paredes ----
$ws.Range("A2").Value = "Pedro"
$ws.Range("B2").Value = "Paredes"
$ws.Range("C2").Value = "paredesp"
$ws.Range("D2").Value = "paredesp@yopmail.com"
$ws.Range("E2").Value = "BC"
$ws.Range("F2").Value = "PAR"
$ws.Range("H2").Value = "BQM"
$ws.Range("I2").Value = "Masculino"
$ws.Range("J2").Value = 45869632
$ws.Range("K2").Value = "Licenciado "
$ws.Range("L2").Value = "Investigación celular"
$ws.Range("N2").Value = 4168546321

# ---- Row 3: Mary James ----
$ws.Range("A3").Value = "Mary"
$ws.Range("B3").Value = "James"
$ws.Range("C3").Value = "mjames"
$ws.Range("D3").Value = "mjames@yopmail.col"
$ws.Range("E3").Value = "BC"
$ws.Range("F3").Value = "PAR"
$ws.Range("H3").Value = "BQM"
$ws.Range("I3").Value = "Femenino"
$ws.Range("J3").Value = 2074967
$ws.Range("K3").Value = "Licenciado"
$ws.Range("L3").Value = "Genética aplicada"
$ws.Range("N3").Value = 4168523651

# ---- Hyperlinks: refresh D2 (email changed) and add the missing D3 one ----
# NB: Range.Hyperlinks.Delete() on a range with no hyperlink wipes the whole
# sheet's hyperlink collection in this runtime, so only call it when a
# hyperlink is actually present on that cell.
if ($ws.Range("D2").Hyperlinks.Count -gt 0) {
    $ws.Range("D2").Hyperlinks.Delete()
}
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:paredesp@yopmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "paredesp@yopmail.com")

if ($ws.Range("D3").Hyperlinks.Count -gt 0) {
    $ws.Range("D3").Hyperlinks.Delete()
}
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:mjames@yopmail.col", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "mjames@yopmail.col")

# Adding a hyperlink re-styles the cell with the built-in "Hyperlink" style
# (blue + underline). The original workbook instead used a plain blue,
# non-underlined font for these cells, so restore that look on both cells.
foreach ($addr in @("D2", "D3")) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 10
    $rng.Font.Color = 16711680
    $rng.Font.Underline = -4142
    $rng.Font.Bold = $false
}

# ---- Column D width ----
$ws.Columns.Item(4).ColumnWidth = 25.166666666666668

# ---- Active selection ----
$ws.Range("F3").Select()
